# "bieu do trinh tu" - add two new Gantt tasks (Sequence Diagram, UI design)
# and correct the "responsible person" assignments on two existing tasks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix "người phụ trách" (person in charge) on existing rows ---------
$ws.Range("D5").Value = "Chung, Hiếu"
$ws.Range("D6").Value = "Hoàn, An"

# --- Row 7: new task "Vẽ Biểu đồ Trình tự" (Draw Sequence Diagram) -----
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Vẽ Biểu đồ Trình tự"
$ws.Range("C7").Value = "Vẽ Biểu đồ Trình tự (Sequence Diagram) để thể hiện luồng thông" + [char]10 + "điệp giữa các đối tượng"
$ws.Range("D7").Value = "An"

# --- Row 8: new task "Thiết kế giao diện" (UI design) ------------------
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Thiết kế giao diện"
$ws.Range("C8").Value = "Phác thảo giao diện cho chức năng"
$ws.Range("D8").Value = "Ninh"

# --- Match formatting to the existing sibling rows ----------------------
$ws.Range("A3").Copy()
$ws.Range("A7:A8").PasteSpecial(-4122)

$ws.Range("B3").Copy()
$ws.Range("B7:B8").PasteSpecial(-4122)

$ws.Range("C4").Copy()
$ws.Range("C7").PasteSpecial(-4122)

$ws.Range("C3").Copy()
$ws.Range("C8").PasteSpecial(-4122)

$ws.Range("D3").Copy()
$ws.Range("D7:D8").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Selection moves to K4, matching the saved file ---------------------
$ws.Range("K4").Select()
